$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'33.800.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +10.11%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.806.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +7.13%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'227.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.33%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.538"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.49%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.15%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'30.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'47.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.19%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.278"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.27%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.68%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.0931"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.63%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'2.066.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.16%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'1.808.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +7.28%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.633"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.19%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'33.810.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +10.14%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'10.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.34%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'4.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.97%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'68.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.97%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'254.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.28%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.0₃0738"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.24%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.09%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'10.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.01%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  +0.15%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'  +0.95%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'158.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.62%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'16.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.12%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +2.73%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.84%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +0.40%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = "'MinaProtocolToken"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'2.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +399.88%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("B32").Value = "'Filecoin"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'3.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +9.31%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  +5.66%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = "'Hedera"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.0506"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.62%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").Value = "'InternetComputer(DFINITY)"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'3.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.11%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "'Maker"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.533.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.52%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'LidoDAOToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.54%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +3.34%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.0185"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.03%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'83.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.35%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.613"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.81%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'2.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.15%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  +0.13%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.37%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  +5.34%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.0519"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.88%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  +4.01%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'1.953.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +7.04%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -0.02%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  +3.82%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'51.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.42%  "
$ws.Range("E51").Style = "Normal"
